$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "展览" — update F-column (想去人数) values
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 1486
$ws.Cells.Item(4, 6).Value = 2114
$ws.Cells.Item(5, 6).Value = 7449
$ws.Cells.Item(7, 6).Value = 4782
$ws.Cells.Item(8, 6).Value = 7034
$ws.Cells.Item(10, 6).Value = 275
$ws.Cells.Item(11, 6).Value = 1489
$ws.Cells.Item(12, 6).Value = 860
$ws.Cells.Item(13, 6).Value = 176
$ws.Cells.Item(17, 6).Value = 161
$ws.Cells.Item(21, 6).Value = 1160
$ws.Cells.Item(23, 6).Value = 3
$ws.Cells.Item(25, 6).Value = 1226
$ws.Cells.Item(27, 6).Value = 143
$ws.Cells.Item(29, 6).Value = 42
$ws.Cells.Item(30, 6).Value = 176
$ws.Cells.Item(32, 6).Value = 40
$ws.Cells.Item(33, 6).Value = 93
$ws.Cells.Item(34, 6).Value = 33
$ws.Cells.Item(35, 6).Value = 549
$ws.Cells.Item(37, 6).Value = 71
$ws.Cells.Item(38, 6).Value = 62
$ws.Cells.Item(39, 6).Value = 374
$ws.Cells.Item(40, 6).Value = 1199
$ws.Cells.Item(41, 6).Value = 576
$ws.Cells.Item(42, 6).Value = 140

# ---------------------------------------------------------------------------
# Sheet 2: "演出" — update F-column values for existing rows, then insert a
# new row 43 for the "梁祝之父" concert (everything from old row 43 onward
# shifts down by one row).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(26, 6).Value = 636
$ws.Cells.Item(31, 6).Value = 855
$ws.Cells.Item(36, 6).Value = 111
$ws.Cells.Item(40, 6).Value = 140
$ws.Cells.Item(42, 6).Value = 10

# Insert the new row before the old row 43 ("四月是你的谎言...").
$ws.Rows.Item(43).Insert()

# The freshly inserted row's A-cell doesn't carry the bordered/bold style
# used by the rest of column A; copy that formatting over explicitly.
$ws.Cells.Item(42, 1).Copy()
$ws.Cells.Item(43, 1).PasteSpecial(-4122)

# Fix up column A (sequential index = row - 1) for the new row and every
# row that shifted down underneath it.
for ($r = 43; $r -le 47; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Populate the newly inserted row 43 with the new concert's data.
$ws.Cells.Item(43, 2).Value = "'2024-07-26"
$ws.Cells.Item(3, 2).Copy()
$ws.Cells.Item(43, 2).PasteSpecial(-4122)

$ws.Cells.Item(43, 3).Value = "上海·梁祝之父：何占豪指挥·《梁祝》65周年大型东方交响音乐会"
$ws.Cells.Item(43, 4).Value = "丁香路425号(上海科技馆地铁站1号口步行460米) 上海东方艺术中心音乐厅"
$ws.Cells.Item(43, 5).Value = "2024.07.26 19:30-07.26 21:00"
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 80
$ws.Cells.Item(43, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85591"
$ws.Cells.Item(43, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/8bMGNbdd1715238003823.jpeg"

# ---------------------------------------------------------------------------
# Sheet 3: "本地生活" — update F-column values
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 6).Value = 730
$ws.Cells.Item(6, 6).Value = 673
$ws.Cells.Item(8, 6).Value = 1592
$ws.Cells.Item(9, 6).Value = 2486

# ---------------------------------------------------------------------------
# Sheet 4: "全部类型" — update F-column values
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 730
$ws.Cells.Item(3, 6).Value = 1486
$ws.Cells.Item(7, 6).Value = 673
$ws.Cells.Item(8, 6).Value = 673
$ws.Cells.Item(9, 6).Value = 7450
$ws.Cells.Item(11, 6).Value = 4782
$ws.Cells.Item(13, 6).Value = 7034
$ws.Cells.Item(14, 6).Value = 275
$ws.Cells.Item(15, 6).Value = 1489
$ws.Cells.Item(16, 6).Value = 860
$ws.Cells.Item(17, 6).Value = 176
$ws.Cells.Item(18, 6).Value = 1592
$ws.Cells.Item(19, 6).Value = 2486
$ws.Cells.Item(20, 6).Value = 204
$ws.Cells.Item(23, 6).Value = 161
$ws.Cells.Item(26, 6).Value = 1160
$ws.Cells.Item(27, 6).Value = 636
$ws.Cells.Item(29, 6).Value = 3
$ws.Cells.Item(30, 6).Value = 1226
$ws.Cells.Item(31, 6).Value = 143
$ws.Cells.Item(32, 6).Value = 176
$ws.Cells.Item(34, 6).Value = 855
$ws.Cells.Item(35, 6).Value = 40
$ws.Cells.Item(36, 6).Value = 93
$ws.Cells.Item(37, 6).Value = 988
$ws.Cells.Item(38, 6).Value = 549
$ws.Cells.Item(40, 6).Value = 71
$ws.Cells.Item(41, 6).Value = 62
$ws.Cells.Item(42, 6).Value = 111
$ws.Cells.Item(43, 6).Value = 374
$ws.Cells.Item(44, 6).Value = 576
$ws.Cells.Item(46, 6).Value = 140
$ws.Cells.Item(48, 6).Value = 140
